# Refined metadata to be additional tab
#
# 1. Re-stamp the per-row "time_taken" column (F) on the "data" sheet with
#    the timestamps captured on the later run of the scraping script.
# 2. Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (name/id/version/timestamps/request URL).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")

# --- 1. Update time_taken timestamps on the "data" sheet -------------------
$ws1.Range("F2").Value  = "2021-10-05 14:22:12.052606"
$ws1.Range("F3").Value  = "2021-10-05 14:22:12.052613"
$ws1.Range("F4").Value  = "2021-10-05 14:22:12.052617"
$ws1.Range("F5").Value  = "2021-10-05 14:22:12.052619"
$ws1.Range("F6").Value  = "2021-10-05 14:22:12.052622"
$ws1.Range("F7").Value  = "2021-10-05 14:22:12.052625"
$ws1.Range("F8").Value  = "2021-10-05 14:22:12.052628"
$ws1.Range("F9").Value  = "2021-10-05 14:22:12.052630"
$ws1.Range("F10").Value = "2021-10-05 14:22:12.052633"
$ws1.Range("F11").Value = "2021-10-05 14:22:12.052635"
$ws1.Range("F12").Value = "2021-10-05 14:22:12.052638"
$ws1.Range("F13").Value = "2021-10-05 14:22:12.052640"
$ws1.Range("F14").Value = "2021-10-05 14:22:12.052643"
$ws1.Range("F15").Value = "2021-10-05 14:22:12.052645"
$ws1.Range("F16").Value = "2021-10-05 14:22:12.052648"
$ws1.Range("F17").Value = "2021-10-05 14:22:12.052650"
$ws1.Range("F18").Value = "2021-10-05 14:22:12.052653"
$ws1.Range("F19").Value = "2021-10-05 14:22:12.052655"
$ws1.Range("F20").Value = "2021-10-05 14:22:12.052658"
$ws1.Range("F21").Value = "2021-10-05 14:22:12.052660"
$ws1.Range("F22").Value = "2021-10-05 14:22:12.052663"
$ws1.Range("F23").Value = "2021-10-05 14:22:12.052665"
$ws1.Range("F24").Value = "2021-10-05 14:22:12.052668"
$ws1.Range("F25").Value = "2021-10-05 14:22:12.052670"
$ws1.Range("F26").Value = "2021-10-05 14:22:12.052673"
$ws1.Range("F27").Value = "2021-10-05 14:22:12.052676"
$ws1.Range("F28").Value = "2021-10-05 14:22:12.052678"
$ws1.Range("F29").Value = "2021-10-05 14:22:12.052681"
$ws1.Range("F30").Value = "2021-10-05 14:22:12.052683"
$ws1.Range("F31").Value = "2021-10-05 14:22:12.052685"
$ws1.Range("F32").Value = "2021-10-05 14:22:12.052688"
$ws1.Range("F33").Value = "2021-10-05 14:22:12.052690"
$ws1.Range("F34").Value = "2021-10-05 14:22:12.052693"
$ws1.Range("F35").Value = "2021-10-05 14:22:12.052696"
$ws1.Range("F36").Value = "2021-10-05 14:22:12.052698"
$ws1.Range("F37").Value = "2021-10-05 14:22:12.052701"
$ws1.Range("F38").Value = "2021-10-05 14:22:12.052703"
$ws1.Range("F39").Value = "2021-10-05 14:22:12.052706"

# --- 2. Add the "metadata" worksheet, right after "data" -------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Copy the header-row look (bold, centered, bordered) from the "data" sheet
# onto the new header row, and the index-column look onto A2.
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)   # xlPasteFormats

$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Peroxisomal disorders"
$ws2.Range("C2").Value = 114

# data_version must stay textual ("1.16"), not be coerced to the number 1.16
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.16"

$ws2.Range("E2").Value = "2021-09-07T14:17:24.235448Z"
$ws2.Range("F2").Value = "2021-10-05 14:22:12.048961"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/114/?format=json"

# Leave the "data" tab selected/active, as in the original workbook.
$ws1.Activate()
